# Tweak spacing/size of logo (slide layout "Title Slide")
#
# EMU -> point conversion used by the PowerPoint object model:
#   1 point = 12700 EMU

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$cl = $s.CustomLayout

# Shape 4: "Title 1" (ctrTitle placeholder)
#   off  5290686,669497  -> 4946904,669497
#   ext 33270141,2560320 -> 33997392,2560320
$shTitle = $cl.Shapes.Item(4)
$shTitle.Left   = 4946904 / 12700
$shTitle.Top    = 669497 / 12700
$shTitle.Width  = 33997392 / 12700
$shTitle.Height = 2560320 / 12700

# Shape 5: "Title 1" (author name text box)
#   off  3383280,3413760 -> 4946904,3413760
#   ext 37307520,1752600 -> 33997392,1752600
$shAuthor = $cl.Shapes.Item(5)
$shAuthor.Left   = 4946904 / 12700
$shAuthor.Top    = 3413760 / 12700
$shAuthor.Width  = 33997392 / 12700
$shAuthor.Height = 1752600 / 12700

# Shape 17: "Picture 6" (university logo picture)
#   off  38875811,309307 -> 39456360,512064
#   ext  4252836,4810539 -> 3920693,4434840
$shPic = $cl.Shapes.Item(17)
$shPic.Left   = 39456360 / 12700
$shPic.Top    = 512064 / 12700
$shPic.Width  = 3920693 / 12700
$shPic.Height = 4434840 / 12700

# Shape 18: "Content Placeholder 33" (additional graphic/logo placeholder)
#   off   762000,356616 -> 512064,512064
#   ext  4252913,4809744 -> 3922776,4434840
$shContent = $cl.Shapes.Item(18)
$shContent.Left   = 512064 / 12700
$shContent.Top    = 512064 / 12700
$shContent.Width  = 3922776 / 12700
$shContent.Height = 4434840 / 12700
